# chore: update Sheets via scheduled runner
# Applies updated currentAveragePrice / LevePrice / LeveProfit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 189.33333
$ws.Range("I2").Value = 209
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 209
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -96
$ws.Range("N2").Value = -376
$ws.Range("H48").Value = 6250
$ws.Range("J48").Value = 6250
$ws.Range("L48").Value = 18750
$ws.Range("N48").Value = -19334
$ws.Range("H51").Value = 1833.3334
$ws.Range("I51").Value = 1666.6666
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 1666.6666
$ws.Range("L51").Value = 2000
$ws.Range("M51").Value = -1182.6666
$ws.Range("N51").Value = -2968
$ws.Range("H52").Value = 10000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 10000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 30000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -30320
$ws.Range("H56").Value = 6250
$ws.Range("J56").Value = 6250
$ws.Range("L56").Value = 18750
$ws.Range("N56").Value = -19818
$ws.Range("H113").Value = 2857.1428
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3200
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -9708
$ws.Range("H121").Value = 1968.3334
$ws.Range("I121").Value = 1000
$ws.Range("J121").Value = 2452.5
$ws.Range("K121").Value = 3000
$ws.Range("L121").Value = 7357.5
$ws.Range("M121").Value = -1253
$ws.Range("N121").Value = -10851.5
$ws.Range("H129").Value = 1177.591
$ws.Range("I129").Value = 767
$ws.Range("J129").Value = 1369.2
$ws.Range("K129").Value = 2301
$ws.Range("L129").Value = 4107.6
$ws.Range("M129").Value = 2699
$ws.Range("N129").Value = -14107.6
$ws.Range("H132").Value = 1954.5264
$ws.Range("I132").Value = 1544.7333
$ws.Range("J132").Value = 3491.25
$ws.Range("K132").Value = 4634.199900000001
$ws.Range("L132").Value = 10473.75
$ws.Range("M132").Value = -2104.199900000001
$ws.Range("N132").Value = -15533.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15485.788
$ws.Range("I32").Value = 16616.191
$ws.Range("J32").Value = 4860
$ws.Range("K32").Value = 16616.191
$ws.Range("L32").Value = 4860
$ws.Range("M32").Value = -16329.191
$ws.Range("N32").Value = -5434
$ws.Range("H44").Value = 39800
$ws.Range("J44").Value = 39800
$ws.Range("L44").Value = 39800
$ws.Range("N44").Value = -40776
$ws.Range("H61").Value = 1569.0526
$ws.Range("I61").Value = 1353.6471
$ws.Range("K61").Value = 1353.6471
$ws.Range("M61").Value = -1141.6471
$ws.Range("H74").Value = 1181.619
$ws.Range("I74").Value = 1095.4736
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1095.4736
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -221.4736
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1181.619
$ws.Range("I77").Value = 1095.4736
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 5477.368
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -1109.368
$ws.Range("N77").Value = -18736
$ws.Range("H122").Value = 1185.4286
$ws.Range("I122").Value = 1185.4286
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3556.2858
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1106.2858
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 30428
$ws.Range("J123").Value = 30428
$ws.Range("L123").Value = 30428
$ws.Range("N123").Value = -40228
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 1569.0526
$ws.Range("I136").Value = 1353.6471
$ws.Range("K136").Value = 4060.9413
$ws.Range("M136").Value = -1510.9413

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2008.1
$ws.Range("I134").Value = 1574.2727
$ws.Range("J134").Value = 3201.125
$ws.Range("K134").Value = 4722.8181
$ws.Range("L134").Value = 9603.375
$ws.Range("M134").Value = -2187.8181
$ws.Range("N134").Value = -14673.375

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1672.8636
$ws.Range("I58").Value = 1727.7222
$ws.Range("J58").Value = 1426
$ws.Range("K58").Value = 1727.7222
$ws.Range("L58").Value = 1426
$ws.Range("M58").Value = -1524.7222
$ws.Range("N58").Value = -1832
$ws.Range("H94").Value = 2818.4443
$ws.Range("I94").Value = 1141.4286
$ws.Range("J94").Value = 3885.6365
$ws.Range("K94").Value = 1141.4286
$ws.Range("L94").Value = 3885.6365
$ws.Range("M94").Value = -690.4286
$ws.Range("N94").Value = -4787.636500000001
$ws.Range("H136").Value = 1672.8636
$ws.Range("I136").Value = 1727.7222
$ws.Range("J136").Value = 1426
$ws.Range("K136").Value = 5183.1666
$ws.Range("L136").Value = 4278
$ws.Range("M136").Value = -2633.1666
$ws.Range("N136").Value = -9378

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 675.5
$ws.Range("I113").Value = 665.5714
$ws.Range("J113").Value = 680.8461
$ws.Range("K113").Value = 1996.7142
$ws.Range("L113").Value = 2042.5383
$ws.Range("M113").Value = 173.2857999999999
$ws.Range("N113").Value = -6382.5383
$ws.Range("H131").Value = 33338414
$ws.Range("I131").Value = 11480
$ws.Range("J131").Value = 47621384
$ws.Range("K131").Value = 34440
$ws.Range("L131").Value = 142864152
$ws.Range("M131").Value = -29400
$ws.Range("N131").Value = -142874232

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 115023.6
$ws.Range("J15").Value = 115023.6
$ws.Range("L15").Value = 115023.6
$ws.Range("N15").Value = -115599.6
$ws.Range("H81").Value = 115023.6
$ws.Range("J81").Value = 115023.6
$ws.Range("L81").Value = 115023.6
$ws.Range("N81").Value = -117019.6
$ws.Range("H84").Value = 115023.6
$ws.Range("J84").Value = 115023.6
$ws.Range("L84").Value = 345070.8
$ws.Range("N84").Value = -355054.8
$ws.Range("H109").Value = 20282.9
$ws.Range("J109").Value = 20282.9
$ws.Range("L109").Value = 20282.9
$ws.Range("N109").Value = -22362.9
$ws.Range("H123").Value = 13103.667
$ws.Range("J123").Value = 13103.667
$ws.Range("L123").Value = 13103.667
$ws.Range("N123").Value = -18003.667

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 235000
$ws.Range("J25").Value = 20000
$ws.Range("L25").Value = 20000
$ws.Range("N25").Value = -20460
$ws.Range("H122").Value = 17864358
$ws.Range("I122").Value = 31256626
$ws.Range("J122").Value = 7999.1665
$ws.Range("K122").Value = 93769878
$ws.Range("L122").Value = 23997.4995
$ws.Range("M122").Value = -93767428
$ws.Range("N122").Value = -28897.4995
$ws.Range("H133").Value = 57531.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 57531.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 57531.332
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -62591.332
$ws.Range("H136").Value = 2608.4333
$ws.Range("I136").Value = 2294.75
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 6884.25
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -4334.25
$ws.Range("N136").Value = -26100

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 50050.625
$ws.Range("J123").Value = 50050.625
$ws.Range("L123").Value = 50050.625
$ws.Range("N123").Value = -59850.625
$ws.Range("H125").Value = 61837
$ws.Range("J125").Value = 61837
$ws.Range("L125").Value = 61837
$ws.Range("N125").Value = -71677
$ws.Range("H132").Value = 1194.909
$ws.Range("I132").Value = 867.8148
$ws.Range("J132").Value = 2666.8333
$ws.Range("K132").Value = 2603.4444
$ws.Range("L132").Value = 8000.499899999999
$ws.Range("M132").Value = -73.44439999999986
$ws.Range("N132").Value = -13060.4999
$ws.Range("H136").Value = 1153.1305
$ws.Range("I136").Value = 1137.3636
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 3412.0908
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -862.0907999999999
$ws.Range("N136").Value = -9600
